# Applies the English -> Italian translation edits described by the diff.
# Each Find/Execute call is scoped to a specific paragraph (or comment) range
# so that only the intended occurrence is replaced.

$d = $word.ActiveDocument

function Replace-InRange($rng, [string]$old, [string]$new) {
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

$paras = $d.Paragraphs

# Paragraph 1: "English / Portuguese / French / Thai / Vietnamese / Spanish"
$p1 = $paras.Item(1).Range
Replace-InRange $p1 "English" "Inglese"
$p1 = $paras.Item(1).Range
Replace-InRange $p1 " / Portuguese / French / Thai / Vietnamese / Spanish" " / Portoghese / Francese / Thai / Vietnamita / Spagnolo"

# Paragraph 3: "English"
$p3 = $paras.Item(3).Range
Replace-InRange $p3 "English" "Inglese"

# Table cell paragraphs (Brief / description / Target audience / description)
$p5 = $paras.Item(5).Range
Replace-InRange $p5 "Brief" "Breve"

$p6 = $paras.Item(6).Range
Replace-InRange $p6 "An email sent to partners in the target country who RSVPed yes but didn’t submit their documents by the deadline. We will be revoking their invites. It will be sent via customer.io" "Un'email inviata ai partner nel paese target che hanno risposto sì ma non hanno inviato i loro documenti entro la scadenza. Revoceremo i loro inviti. Sarà inviato tramite customer.io"

$p8 = $paras.Item(8).Range
Replace-InRange $p8 "Target audience" "Pubblico target"

$p9 = $paras.Item(9).Range
Replace-InRange $p9 "Invited partners who didn’t submit their documents on time" "Partner invitati che non hanno inviato i loro documenti in tempo"

# Paragraph 12: "Subject line: Your [EVENT NAME] registration"
#   -> "Oggetto: La tua registrazione per [EVENT NAME]"
$p12 = $paras.Item(12).Range
Replace-InRange $p12 "Subject line" "Oggetto"
$p12 = $paras.Item(12).Range
Replace-InRange $p12 ": Your " ": La tua registrazione per "
$p12 = $paras.Item(12).Range
Replace-InRange $p12 " registration" ""

# Paragraph 14: "We didn’t receive your documents on time"
$p14 = $paras.Item(14).Range
Replace-InRange $p14 "We didn’t receive your documents on time" "Non abbiamo ricevuto i tuoi documenti in tempo"

# Paragraph 16: "Hi "
$p16 = $paras.Item(16).Range
Replace-InRange $p16 "Hi " "Ciao "

# Paragraph 17
$p17 = $paras.Item(17).Range
Replace-InRange $p17 "We didn’t receive your documents by the deadline (" "Non abbiamo ricevuto i tuoi documenti entro la scadenza ("
$p17 = $paras.Item(17).Range
Replace-InRange $p17 "). Unfortunately, we’re unable to proceed with your registration for the " "). Purtroppo non possiamo procedere con la tua registrazione per "

# Paragraph 18
$p18 = $paras.Item(18).Range
Replace-InRange $p18 "We wish you the best and hope to see you at our next " "Ti facciamo un grosso in bocca al lupo e speriamo di vederti alla nostra prossima "
$p18 = $paras.Item(18).Range
Replace-InRange $p18 "conference/seminar/affiliate trip" "conferenza/seminario/viaggio per associati"

# Paragraph 19
$p19 = $paras.Item(19).Range
Replace-InRange $p19 "If you have any questions, please contact us via " "Se hai domande, non esitare a contattarci tramite "
$p19 = $paras.Item(19).Range
Replace-InRange $p19 "live chat" "chat live"
$p19 = $paras.Item(19).Range
Replace-InRange $p19 " or " " o "

# Paragraph 20
$p20 = $paras.Item(20).Range
Replace-InRange $p20 "If you have any questions, please contact your country manager, " "In caso di domande, contatta il tuo country manager, "
$p20 = $paras.Item(20).Range
Replace-InRange $p20 ", at " ", all'indirizzo "
$p20 = $paras.Item(20).Range
Replace-InRange $p20 " or " " o al numero "

# Comments: "choose either one" -> "scegli uno dei due" (both comments)
$comments = $d.Comments
for ($i = 1; $i -le $comments.Count; $i++) {
    $cRng = $comments.Item($i).Range
    Replace-InRange $cRng "choose either one" "scegli uno dei due"
}
